$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "A" + " " + "slide" (3 runs) -> single run "A slide"
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "temp"
$title.Text = "A slide"

# "Just" + " " + "an" + " " + "image" + " " + "on" + " " + "this" + " " + "side" (11 runs) -> single run
$caption = $s.Shapes.Item(4).TextFrame.TextRange
$caption.Text = "temp"
$caption.Text = "Just an image on this side"
